# The workbook has a "Table1" table (A1:D6) with columns Date, Song,
# Credit, Topic. This edit removes the "Credit" column entirely: the
# column is deleted from the sheet, the table shrinks to A1:C6 (Date,
# Song, Topic), and the now-unused "Credit"/"HFBC Hymnal Book" shared
# strings disappear because nothing references them anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Remove the whole "Credit" column (column C) from the worksheet; this
# shifts Topic (old column D) left into column C and drops the stale
# trailing column.
$ws.Columns.Item(3).Delete()

# The table still thinks it spans 4 columns (Date, Song, Credit, Topic)
# even though the sheet only has 3 columns of data now - shrink it back
# down to match the new data extent.
$tbl.Resize($ws.Range("A1:C6"))

# Re-stamp the header text for the (former "Credit", now showing
# "Topic") third column so the table's column metadata/name tracks the
# actual header cell text again.
$tbl.HeaderRowRange.Item(1, 3).Value = "Topic"

# Park the selection where the author's last save left it.
$ws.Range("E18").Select()
